$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 34: date, note text, hours
$ws.Range("B34").Value = (Get-Date -Year 2024 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C34").Value = "bugfix in Časové údaje"
$ws.Range("E34").Value = 1

# Move the active selection to E35, matching the recorded cursor position after edit
$ws.Range("E35").Select()
